# Weekly refresh: prepend two new price records (rows 320-321) for
# "Feria Lagunitas de Puerto Montt - Repollo", pushing all later rows
# down by two (old row 320 -> new row 322, etc.) and extending the
# sheet's used range from A1:R415 to A1:R417.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 320.
$ws.Rows.Item(320).Insert()
$ws.Rows.Item(320).Insert()

# New row 320: Copenhague / Primera
$ws.Cells.Item(320, 1).Value2  = 4
$ws.Cells.Item(320, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(320, 3).Value2  = "Los Lagos"
$ws.Cells.Item(320, 4).Value2  = 44663
$ws.Cells.Item(320, 5).Value2  = 10
$ws.Cells.Item(320, 6).Value2  = 100112006
$ws.Cells.Item(320, 7).Value2  = "Repollo"
$ws.Cells.Item(320, 8).Value2  = "Copenhague"
$ws.Cells.Item(320, 9).Value2  = "Primera"
$ws.Cells.Item(320, 10).Value2 = 500
$ws.Cells.Item(320, 11).Value2 = 1900
$ws.Cells.Item(320, 12).Value2 = 1900
$ws.Cells.Item(320, 13).Value2 = 1900
$ws.Cells.Item(320, 14).Value2 = "`$/unidad"
$ws.Cells.Item(320, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(320, 16).Value2 = 1900
$ws.Cells.Item(320, 17).Value2 = 1
$ws.Cells.Item(320, 18).Value2 = "Hortaliza"

# New row 321: Crespo record / Primera
$ws.Cells.Item(321, 1).Value2  = 4
$ws.Cells.Item(321, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(321, 3).Value2  = "Los Lagos"
$ws.Cells.Item(321, 4).Value2  = 44663
$ws.Cells.Item(321, 5).Value2  = 10
$ws.Cells.Item(321, 6).Value2  = 100112006
$ws.Cells.Item(321, 7).Value2  = "Repollo"
$ws.Cells.Item(321, 8).Value2  = "Crespo record"
$ws.Cells.Item(321, 9).Value2  = "Primera"
$ws.Cells.Item(321, 10).Value2 = 800
$ws.Cells.Item(321, 11).Value2 = 1700
$ws.Cells.Item(321, 12).Value2 = 1800
$ws.Cells.Item(321, 13).Value2 = 1750
$ws.Cells.Item(321, 14).Value2 = "`$/unidad"
$ws.Cells.Item(321, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(321, 16).Value2 = 1750
$ws.Cells.Item(321, 17).Value2 = 1
$ws.Cells.Item(321, 18).Value2 = "Hortaliza"
